# Fruta / hortaliza, semanal
# Update the weekly price records: each row's Fecha/Volumen/Precio
# minimo/maximo/promedio/Origen/Precio-$-Kg values are re-shuffled
# to the new weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2..11), keyed by column letter:
#   D = Fecha (serial date number)
#   M = Volumen
#   N = Precio minimo
#   O = Precio maximo
#   P = Precio promedio ponderado
#   R = Origen
#   S = Precio $/Kg
$rows = @{
    2  = @{ D = 44208; M = 85;  N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 }
    3  = @{ D = 44188; M = 150; N = 3000; O = 3400; P = 3240; R = "Provincia de Linares"; S = 1620 }
    4  = @{ D = 44232; M = 200; N = 3000; O = 3000; P = 3000; R = "Provincia de Curicó";  S = 1500 }
    5  = @{ D = 44231; M = 150; N = 3400; O = 3400; P = 3400; R = "Provincia de Curicó";  S = 1700 }
    6  = @{ D = 44174; M = 200; N = 3200; O = 3200; P = 3200; R = "Provincia de Curicó";  S = 1600 }
    7  = @{ D = 44236; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó";  S = 1900 }
    8  = @{ D = 44168; M = 170; N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 }
    9  = @{ D = 44194; M = 120; N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 }
    10 = @{ D = 44237; M = 100; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó";  S = 1900 }
    11 = @{ D = 44238; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó";  S = 1900 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value2  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 13).Value2 = $vals.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value2 = $vals.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value2 = $vals.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value2 = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 18).Value  = $vals.R   # R: Origen
    $ws.Cells.Item($r, 19).Value2 = $vals.S   # S: Precio $/Kg
}
